$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "***maa://25695 (19.41), **maa://32237 (38.89), ***maa://34206 (14.29), ***maa://39951 (19.05), ***maa://39243 (25.0)"
$ws.Range("S11").Value = "maa://22747 (95.0), maa://22501 (98.08)"
$ws.Range("W11").Value = "maa://36713 (97.8)"
$ws.Range("C12").Value = "maa://30766 (88.89), **maa://36678 (50.0)"
$ws.Range("AA12").Value = "maa://23669 (95.83), maa://36677 (94.74), maa://39872 (83.33)"
$ws.Range("G13").Value = "*maa://21248 (75.74), **maa://22728 (47.62)"
$ws.Range("W13").Value = "*maa://34957 (78.05), *maa://22768 (53.33)"
$ws.Range("C14").Value = "maa://30764 (85.37)"
$ws.Range("S16").Value = "maa://22729 (95.14), *maa://28648 (69.09), *maa://36674 (76.92)"
$ws.Range("C17").Value = "maa://21624 (80.65)"
$ws.Range("G17").Value = "maa://22430 (88.57), maa://39599 (82.35)"
$ws.Range("G18").Value = "maa://24421 (90.43)"
$ws.Range("AA19").Value = "*maa://30709 (60.38), *maa://36668 (52.17)"
$ws.Range("K20").Value = "maa://41331 (88.0)"
$ws.Range("O20").Value = "maa://37442 (96.43)"
$ws.Range("AA21").Value = "*maa://21443 (78.83), ***maa://23820 (29.63)"
$ws.Range("AE21").Value = "maa://22524 (94.25), *maa://22432 (74.07)"
$ws.Range("K23").Value = "maa://39756 (91.67), maa://39875 (95.45)"
$ws.Range("O23").Value = "maa://30587 (91.62), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (77.78)"
$ws.Range("C24").Value = "maa://24368 (80.5)"
$ws.Range("C28").Value = "maa://24465 (90.32), maa://25725 (82.05)"
$ws.Range("K28").Value = "*maa://30770 (78.05)"
$ws.Range("W28").Value = "maa://39929 (85.96), ***maa://39723 (15.15), *maa://41749 (66.67)"
$ws.Range("AE28").Value = "maa://36660 (94.07), *maa://36701 (64.0)"
$ws.Range("C29").Value = "maa://31694 (97.87)"
$ws.Range("K29").Value = "maa://28432 (93.4), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)"
$ws.Range("K30").Value = "maa://30442 (94.23)"
$ws.Range("K31").Value = "maa://35926 (93.48), maa://36258 (80.26)"
$ws.Range("S32").Value = "maa://41108 (89.47), maa://41238 (94.12)"
$ws.Range("K35").Value = "maa://41296 (96.77)"
$ws.Range("G36").Value = "maa://24375 (92.31)"
$ws.Range("G39").Value = "maa://25199 (86.11), maa://36670 (88.71), maa://30434 (88.89), ***maa://25036 (16.0)"
$ws.Range("G41").Value = "maa://24466 (94.87)"
$ws.Range("O45").Value = "*maa://36237 (54.55)"
$ws.Range("G46").Value = "maa://35931 (92.41)"
$ws.Range("G51").Value = "*maa://30769 (78.57)"
$ws.Range("G53").Value = "maa://32534 (93.18), **maa://32434 (36.36)"
$ws.Range("G55").Value = "maa://32532 (92.24)"
$ws.Range("G58").Value = "*maa://37964 (61.11)"
$ws.Range("G60").Value = "**maa://40438 (45.45)"

# Row 26 special case: B26 changes from "0" to "1" (kept as text, like its
# sibling count cells), and C26 gains new content.
# A plain .Value = "1" assignment gets auto-coerced to a Number by the
# numeric-literal sniffer, so instead copy the already-text "1" from F26
# (same row, same style) and paste only the value into B26 - this preserves
# both the original style and the text cell type.
$ws.Range("F26").Copy()
$ws.Range("B26").PasteSpecial(-4163)
$ws.Range("C26").Value = "maa://41802 (100.0)"
